$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.445647641019636
$ws.Cells.Item(2, 3).Value = 0.3048912486333797
$ws.Cells.Item(2, 4).Value = 18.71679738969934
$ws.Cells.Item(2, 5).Value = 13.86384647080068
$ws.Cells.Item(2, 7).Value = 34.33118275015303

$ws.Cells.Item(3, 2).Value = 1.445647641019636
$ws.Cells.Item(3, 3).Value = 1.626987699542094
$ws.Cells.Item(3, 4).Value = 3.223369029078222
$ws.Cells.Item(3, 5).Value = 0.5333859586016987
$ws.Cells.Item(3, 7).Value = 6.82939032824165

$ws.Cells.Item(4, 2).Value = 3.272327238179451
$ws.Cells.Item(4, 3).Value = 1.626987699542094
$ws.Cells.Item(4, 4).Value = 18.71679738969934
$ws.Cells.Item(4, 5).Value = 13.86384647080068
$ws.Cells.Item(4, 7).Value = 37.47995879822157

$ws.Cells.Item(5, 2).Value = 3.272327238179451
$ws.Cells.Item(5, 3).Value = 1.626987699542094
$ws.Cells.Item(5, 4).Value = 3.223369029078222
$ws.Cells.Item(5, 5).Value = 0.5333859586016987
$ws.Cells.Item(5, 7).Value = 8.656069925401464

$ws.Cells.Item(6, 2).Value = 3.272327238179451
$ws.Cells.Item(6, 3).Value = 1.626987699542094
$ws.Cells.Item(6, 4).Value = 0.7210945179870265
$ws.Cells.Item(6, 5).Value = 0.5333859586016987
$ws.Cells.Item(6, 7).Value = 6.15379541431027

$ws.Cells.Item(7, 2).Value = 3.272327238179451
$ws.Cells.Item(7, 3).Value = 1.626987699542094
$ws.Cells.Item(7, 4).Value = 3.223369029078222
$ws.Cells.Item(7, 5).Value = 0.5333859586016987
$ws.Cells.Item(7, 7).Value = 8.656069925401464

$ws.Cells.Item(8, 2).Value = 3.272327238179451
$ws.Cells.Item(8, 3).Value = 1.626987699542094
$ws.Cells.Item(8, 4).Value = 0.7210945179870265
$ws.Cells.Item(8, 5).Value = 0.5333859586016987
$ws.Cells.Item(8, 7).Value = 6.15379541431027

$ws.Cells.Item(9, 2).Value = 3.272327238179451
$ws.Cells.Item(9, 3).Value = 1.626987699542094
$ws.Cells.Item(9, 4).Value = 0.7210945179870265
$ws.Cells.Item(9, 5).Value = 0.5333859586016987
$ws.Cells.Item(9, 7).Value = 6.15379541431027

$ws.Cells.Item(10, 2).Value = 3.272327238179451
$ws.Cells.Item(10, 3).Value = 1.626987699542094
$ws.Cells.Item(10, 4).Value = 0.7210945179870265
$ws.Cells.Item(10, 5).Value = 0.5333859586016987
$ws.Cells.Item(10, 7).Value = 6.15379541431027

$ws.Cells.Item(11, 2).Value = 0.1169995834814548
$ws.Cells.Item(11, 3).Value = 0.04103571897497393
$ws.Cells.Item(11, 4).Value = 3.223369029078222
$ws.Cells.Item(11, 5).Value = 0.5333859586016987
$ws.Cells.Item(11, 7).Value = 3.914790290136349

$ws.Cells.Item(12, 2).Value = 0.04172184405617529
$ws.Cells.Item(12, 3).Value = 0.04103571897497393
$ws.Cells.Item(12, 4).Value = 3.223369029078222
$ws.Cells.Item(12, 5).Value = 0.5333859586016987
$ws.Cells.Item(12, 7).Value = 3.83951255071107

$ws.Cells.Item(13, 2).Value = 3.272327238179451
$ws.Cells.Item(13, 3).Value = 1.626987699542094
$ws.Cells.Item(13, 4).Value = 3.223369029078222
$ws.Cells.Item(13, 5).Value = 0.5333859586016987
$ws.Cells.Item(13, 7).Value = 8.656069925401464

$ws.Cells.Item(14, 2).Value = 3.272327238179451
$ws.Cells.Item(14, 3).Value = 1.626987699542094
$ws.Cells.Item(14, 4).Value = 0.7210945179870265
$ws.Cells.Item(14, 5).Value = 0.5333859586016987
$ws.Cells.Item(14, 7).Value = 6.15379541431027

$ws.Cells.Item(15, 2).Value = 1.445647641019636
$ws.Cells.Item(15, 3).Value = 1.626987699542094
$ws.Cells.Item(15, 4).Value = 0.1496068669990043
$ws.Cells.Item(15, 5).Value = 0.5333859586016987
$ws.Cells.Item(15, 7).Value = 3.755628166162433

$ws.Cells.Item(16, 2).Value = 3.272327238179451
$ws.Cells.Item(16, 3).Value = 1.626987699542094
$ws.Cells.Item(16, 4).Value = 18.71679738969934
$ws.Cells.Item(16, 5).Value = 0.5333859586016987
$ws.Cells.Item(16, 7).Value = 24.14949828602258

$ws.Cells.Item(17, 2).Value = 1.445647641019636
$ws.Cells.Item(17, 3).Value = 1.626987699542094
$ws.Cells.Item(17, 4).Value = 18.71679738969934
$ws.Cells.Item(17, 5).Value = 0.5333859586016987
$ws.Cells.Item(17, 7).Value = 22.32281868886277

$ws.Cells.Item(18, 2).Value = 0.6545652718822623
$ws.Cells.Item(18, 3).Value = 1.626987699542094
$ws.Cells.Item(18, 4).Value = 0.7210945179870265
$ws.Cells.Item(18, 5).Value = 0.5333859586016987
$ws.Cells.Item(18, 7).Value = 3.536033448013082

$ws.Cells.Item(19, 2).Value = 0.6545652718822623
$ws.Cells.Item(19, 3).Value = 1.626987699542094
$ws.Cells.Item(19, 4).Value = 3.223369029078222
$ws.Cells.Item(19, 5).Value = 0.5333859586016987
$ws.Cells.Item(19, 7).Value = 6.038307959104277

$ws.Cells.Item(20, 2).Value = 3.272327238179451
$ws.Cells.Item(20, 3).Value = 1.626987699542094
$ws.Cells.Item(20, 4).Value = 3.223369029078222
$ws.Cells.Item(20, 5).Value = 0.5333859586016987
$ws.Cells.Item(20, 7).Value = 8.656069925401464

$ws.Cells.Item(21, 2).Value = 0.00009552326474482342
$ws.Cells.Item(21, 3).Value = 0.04103571897497393
$ws.Cells.Item(21, 4).Value = 18.71679738969934
$ws.Cells.Item(21, 5).Value = 0.5333859586016987
$ws.Cells.Item(21, 7).Value = 19.29131459054075

$ws.Cells.Item(22, 2).Value = 3.272327238179451
$ws.Cells.Item(22, 3).Value = 1.626987699542094
$ws.Cells.Item(22, 4).Value = 0.7210945179870265
$ws.Cells.Item(22, 5).Value = 0.5333859586016987
$ws.Cells.Item(22, 7).Value = 6.15379541431027

$ws.Cells.Item(23, 2).Value = 1.445647641019636
$ws.Cells.Item(23, 3).Value = 0.3048912486333797
$ws.Cells.Item(23, 4).Value = 0.7210945179870265
$ws.Cells.Item(23, 5).Value = 0.5333859586016987
$ws.Cells.Item(23, 7).Value = 3.005019366241741

$ws.Cells.Item(24, 2).Value = 0.04172184405617529
$ws.Cells.Item(24, 3).Value = 0.3048912486333797
$ws.Cells.Item(24, 4).Value = 0.1496068669990043
$ws.Cells.Item(24, 5).Value = 0.5333859586016987
$ws.Cells.Item(24, 7).Value = 1.029605918290258

$ws.Cells.Item(25, 2).Value = 1.445647641019636
$ws.Cells.Item(25, 3).Value = 0.04103571897497393
$ws.Cells.Item(25, 4).Value = 0.1496068669990043
$ws.Cells.Item(25, 5).Value = 0.5333859586016987
$ws.Cells.Item(25, 7).Value = 2.169676185595313

$ws.Cells.Item(26, 2).Value = 1.445647641019636
$ws.Cells.Item(26, 3).Value = 1.626987699542094
$ws.Cells.Item(26, 4).Value = 0.7210945179870265
$ws.Cells.Item(26, 5).Value = 0.5333859586016987
$ws.Cells.Item(26, 7).Value = 4.327115817150455

$ws.Cells.Item(27, 2).Value = 3.272327238179451
$ws.Cells.Item(27, 3).Value = 1.626987699542094
$ws.Cells.Item(27, 4).Value = 0.7210945179870265
$ws.Cells.Item(27, 5).Value = 0.5333859586016987
$ws.Cells.Item(27, 7).Value = 6.15379541431027

$ws.Cells.Item(28, 2).Value = 0.6545652718822623
$ws.Cells.Item(28, 3).Value = 1.626987699542094
$ws.Cells.Item(28, 4).Value = 3.223369029078222
$ws.Cells.Item(28, 5).Value = 13.86384647080068
$ws.Cells.Item(28, 7).Value = 19.36876847130326

$ws.Cells.Item(29, 2).Value = 0.6545652718822623
$ws.Cells.Item(29, 3).Value = 1.626987699542094
$ws.Cells.Item(29, 4).Value = 0.7210945179870265
$ws.Cells.Item(29, 5).Value = 0.5333859586016987
$ws.Cells.Item(29, 7).Value = 3.536033448013082

$ws.Cells.Item(30, 2).Value = 1.445647641019636
$ws.Cells.Item(30, 3).Value = 1.626987699542094
$ws.Cells.Item(30, 4).Value = 0.7210945179870265
$ws.Cells.Item(30, 5).Value = 0.5333859586016987
$ws.Cells.Item(30, 7).Value = 4.327115817150455

$ws.Cells.Item(31, 2).Value = 0.6545652718822623
$ws.Cells.Item(31, 3).Value = 22898927661.19195
$ws.Cells.Item(31, 4).Value = 18.71679738969934
$ws.Cells.Item(31, 5).Value = 13.86384647080068
$ws.Cells.Item(31, 7).Value = 22898927694.42716

$ws.Cells.Item(32, 2).Value = 3.272327238179451
$ws.Cells.Item(32, 3).Value = 1.626987699542094
$ws.Cells.Item(32, 4).Value = 0.7210945179870265
$ws.Cells.Item(32, 5).Value = 0.5333859586016987
$ws.Cells.Item(32, 7).Value = 6.15379541431027

$ws.Cells.Item(33, 2).Value = 1.445647641019636
$ws.Cells.Item(33, 3).Value = 1.626987699542094
$ws.Cells.Item(33, 4).Value = 0.7210945179870265
$ws.Cells.Item(33, 5).Value = 0.5333859586016987
$ws.Cells.Item(33, 7).Value = 4.327115817150455

$ws.Cells.Item(34, 2).Value = 0.6545652718822623
$ws.Cells.Item(34, 3).Value = 0.3048912486333797
$ws.Cells.Item(34, 4).Value = 189.6080260415259
$ws.Cells.Item(34, 5).Value = 0.5333859586016987
$ws.Cells.Item(34, 7).Value = 191.1008685206432

$ws.Cells.Item(35, 2).Value = 3.272327238179451
$ws.Cells.Item(35, 3).Value = 1.626987699542094
$ws.Cells.Item(35, 4).Value = 18.71679738969934
$ws.Cells.Item(35, 5).Value = 0.5333859586016987
$ws.Cells.Item(35, 7).Value = 24.14949828602258

$ws.Cells.Item(36, 2).Value = 3.272327238179451
$ws.Cells.Item(36, 3).Value = 9.983522426115931
$ws.Cells.Item(36, 4).Value = 0.7210945179870265
$ws.Cells.Item(36, 5).Value = 13.86384647080068
$ws.Cells.Item(36, 7).Value = 27.84079065308309

$ws.Cells.Item(37, 2).Value = 1.445647641019636
$ws.Cells.Item(37, 3).Value = 1.626987699542094
$ws.Cells.Item(37, 4).Value = 3.223369029078222
$ws.Cells.Item(37, 5).Value = 0.5333859586016987
$ws.Cells.Item(37, 7).Value = 6.82939032824165

$ws.Cells.Item(38, 2).Value = 3.272327238179451
$ws.Cells.Item(38, 3).Value = 1.626987699542094
$ws.Cells.Item(38, 4).Value = 3993.344853322108
$ws.Cells.Item(38, 5).Value = 0.5333859586016987
$ws.Cells.Item(38, 7).Value = 3998.777554218431

$ws.Cells.Item(39, 2).Value = 1.445647641019636
$ws.Cells.Item(39, 3).Value = 1.626987699542094
$ws.Cells.Item(39, 4).Value = 0.7210945179870265
$ws.Cells.Item(39, 5).Value = 0.5333859586016987
$ws.Cells.Item(39, 7).Value = 4.327115817150455

$ws.Cells.Item(40, 2).Value = 1.445647641019636
$ws.Cells.Item(40, 3).Value = 1.626987699542094
$ws.Cells.Item(40, 4).Value = 0.7210945179870265
$ws.Cells.Item(40, 5).Value = 0.5333859586016987
$ws.Cells.Item(40, 7).Value = 4.327115817150455

$ws.Cells.Item(41, 2).Value = 3.272327238179451
$ws.Cells.Item(41, 3).Value = 1.626987699542094
$ws.Cells.Item(41, 4).Value = 3.223369029078222
$ws.Cells.Item(41, 5).Value = 0.5333859586016987
$ws.Cells.Item(41, 7).Value = 8.656069925401464

$ws.Cells.Item(42, 2).Value = 1.445647641019636
$ws.Cells.Item(42, 3).Value = 1.626987699542094
$ws.Cells.Item(42, 4).Value = 18.71679738969934
$ws.Cells.Item(42, 5).Value = 0.5333859586016987
$ws.Cells.Item(42, 7).Value = 22.32281868886277
